# Typeahead search (in progress...)
# Append a new data row (2015.11.21 / 2) below the existing "Date"/"Exercises"
# header row on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A's new cell to be stored as literal text instead of having
# Excel auto-convert the dotted date-looking string into a date serial
# number: pre-format as Text, assign the value, then drop the now-unneeded
# cell formatting so the new row keeps the sheet's default (unstyled) look.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "2015.11.21"
$ws.Range("A2").ClearFormats()

$ws.Range("B2").Value = 2

# Move the active selection down onto the newly added row, like a user
# would after typing the entry.
$null = $ws.Range("A2").Select()
